# Update "Benchmark AnyLogic API" workbook:
#  - add a new "responsetijd private cloud" column (G) of data next to the
#    existing benchmark columns
#  - plot that new column as a third series (with trendline) on the existing
#    scatter chart
#  - enlarge the chart so the extra series/legend entry still fits

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. worksheet data -----------------------------------------------------
# Header + values for the new "responsetijd private cloud" series.
$ws.Range("G2").Value = "responsetijd private cloud"
$ws.Range("G3").Value = 5.4
$ws.Range("G4").Value = 6.2
$ws.Range("G5").Value = 12.8

# Leave the last-used selection on the newly added cell, like the source file.
[void]$ws.Range("G5").Select()

# --- 2. chart: new series ----------------------------------------------------
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$seriesCollection = $chart.SeriesCollection()

$newSeries = $seriesCollection.NewSeries()
$newSeries.Formula = "=SERIES(Sheet1!`$G`$2,Sheet1!`$C`$3:`$C`$5,Sheet1!`$G`$3:`$G`$5,3)"

# Match the other two series, which each carry a linear trendline.
$trend = $newSeries.Trendlines().Add(-4132)
$trend.Type = -4132

# --- 3. resize the chart so the third series/legend entry has room ----------
# Anchored "to" cell moves from col J (9) / 304800 EMU to col N (12) / 60960 EMU.
$co.Width = 776.83515625
